$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E (new quarter-end data is reported first,
# pushing the existing quarterly columns two places to the right).
$ws.Columns("D:E").Insert()

# Copy the number formats (date format for the header row, the "#,##0"
# style for the data rows, etc.) from column F - which now holds what used
# to be column D - onto the two freshly inserted blank columns.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New values for the two inserted columns (D = most recent quarter,
# E = the quarter before it), keyed by row number.
$newData = @{
    7 = @(43496, 43404)
    8 = @(1284000, 1294000)
    9 = @(572000, 587000)
    10 = @(712000, 707000)
    12 = @(102000, 104000)
    13 = @(0, 0)
    14 = @("NA", 0)
    15 = @(0, 0)
    17 = @(1034000, 1045000)
    18 = @(250000, 249000)
    20 = @(16000, 13000)
    21 = @(320000, 318000)
    22 = @(18000, 18000)
    23 = @(248000, 244000)
    24 = @(-256000, 30000)
    25 = @(0, 0)
    26 = @(504000, 214000)
    27 = @(504000, 214000)
    28 = @(0, 0)
    29 = @("NA", -19000)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-16000, -13000)
    33 = @(504000, 195000)
    34 = @(0, 0)
    35 = @(504000, 195000)
    38 = @(43496, 43404)
    41 = @(2057000, 2247000)
    42 = @(0, 0)
    43 = @(833000, 776000)
    44 = @(653000, 638000)
    45 = @(169000, 187000)
    46 = @(3712000, 3848000)
    47 = @(77000, 68000)
    48 = @(829000, 822000)
    49 = @(3699000, 3464000)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(635000, 339000)
    53 = @(0, 0)
    54 = @(8952000, 8541000)
    57 = @(315000, 340000)
    58 = @("NA", 0)
    59 = @(780000, 831000)
    60 = @(1095000, 1171000)
    61 = @(1798000, 1799000)
    62 = @(1023000, 1000000)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(3916000, 3974000)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(90000, -336000)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(5036000, 4567000)
    77 = @(0, 0)
    80 = @(43496, 43404)
    81 = @(504000, 195000)
    83 = @(54000, 56000)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(213000, 372000)
    91 = @(-39000, -36000)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-290000, -114000)
    96 = @(-52000, -47000)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-122000, -131000)
    101 = @(9000, -11000)
    102 = @(-190000, 116000)

}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
}
